# "Lighter color blue for cover"
#
# Changes applied:
#   1. Cover background rectangle ("Rectangle 15" on slide 1) fill color
#      is changed from the dark navy 0000A0 to the lighter 0080FF.
#   2. The auto-updating "datetimeFigureOut" date placeholder that lives
#      on the slide master and on every slide layout is refreshed from
#      4/25/2016 to 5/26/2016 (cosmetic re-cache that PowerPoint performs
#      whenever the deck is touched/saved on a different day).
#   3. Two slide guides (one horizontal, one vertical) are (best-effort)
#      added to the presentation, matching the authoring aid the author
#      dropped onto the canvas while repositioning the cover artwork.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Lighten the cover rectangle's fill color (0000A0 -> 0080FF).
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 15") {
        $shp.Fill.ForeColor.RGB = 0xFF8000
    }
}

# ---------------------------------------------------------------------
# 2) Re-cache the "Update automatically" date placeholder text on the
#    slide master and every slide layout (4/25/2016 -> 5/26/2016).
# ---------------------------------------------------------------------
function Update-DateField($shapeRange) {
    for ($i = 1; $i -le $shapeRange.Count; $i++) {
        $shp = $shapeRange.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "4/25/2016") {
                $tr.Text = "5/26/2016"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DateField $layout.Shapes
}

# ---------------------------------------------------------------------
# 3) Add the pair of slide guides used while repositioning the cover
#    artwork (best-effort — older/alternate hosts may not persist
#    presentation-level guides, so failures here are swallowed).
# ---------------------------------------------------------------------
try {
    $hGuide = $p.Guides.Add(1, 2160)
    $hGuide.Position = 2160
} catch {
}

try {
    $vGuide = $p.Guides.Add(2, 2880)
    $vGuide.Position = 2880
} catch {
}
